# Convert EMU to points for the PowerPoint COM object model (Left/Top/Width/Height
# are expressed in points = EMU / 12700). A tiny half-EMU epsilon is added before
# the division so that the point value, once it round-trips through the host's
# internal float representation, lands back on the exact target EMU instead of
# being truncated down by one unit.
function Emu-ToPoints($emu) {
    return ($emu + 0.5) / 12700
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape "Rectangle 65" - large background rounded rectangle widened/shifted left
$rect = $s.Shapes.Item("Rectangle 65")
$rect.Left   = Emu-ToPoints 563878
$rect.Top    = Emu-ToPoints 480098
$rect.Width  = Emu-ToPoints 9157933
$rect.Height = Emu-ToPoints 5988339

# Shape "Straight Arrow Connector 15" - connector moved/resized to match
$conn = $s.Shapes.Item("Straight Arrow Connector 15")
$conn.Left   = Emu-ToPoints 563878
$conn.Top    = Emu-ToPoints 6103658
$conn.Width  = Emu-ToPoints 2659120
$conn.Height = Emu-ToPoints 13531
